$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B2: remove trailing semicolon from the email address
$ws.Range('B2').Value = 'emailbit21@gmail.com'

# Copy formatting (style s=2, Arial 10 left-aligned) from row 2 down to the new rows 3:36
$ws.Range('A2:B2').Copy()
$ws.Range('A3:B36').PasteSpecial(-4122)

# Populate the new company rows (site, email)
$ws.Range('A3').Value = 'https://condolivre.com.br/gestao-financeira/quanto-cobrar-por-manutencao-de-elevadores/#:~:text=A%20manuten%C3%A7%C3%A3o%20preventiva%20de%20elevadores,e%20da%20complexidade%20do%20equipamento.'
$ws.Range('B3').Value = 'contato@condolivre.com.br'
$ws.Range('A4').Value = 'https://coteibem.sindiconet.com.br/fornecedores/manutencao-elevadores/sp/sao-paulo'
$ws.Range('B4').Value = 'contato@coteibem.com.br'
$ws.Range('A5').Value = 'https://spelevadores.com.br/'
$ws.Range('B5').Value = 'contato@spelevadores.com.br'
$ws.Range('A6').Value = 'https://villarta.com.br/elevadores-e-escadas-rolantes-villarta/lista-de-empresas-de-elevadores-em-sp/'
$ws.Range('B6').Value = 'protecaodedados@villarta.com.br'
$ws.Range('A7').Value = 'https://www.basselevadores.com.br/elevadores-sao-paulo-sp.php'
$ws.Range('B7').Value = 'contato@basselevadores.com.br'
$ws.Range('A8').Value = 'https://elevadoresoiwa.com.br/'
$ws.Range('B8').Value = 'oiwa@elevadoresoiwa.com.br;comercial@elevadoresoiwa.com.br'
$ws.Range('A9').Value = 'https://www.otis.com/pt/br'
$ws.Range('B9').Value = 'cac@otis.com;navigati_cac@otis.com;imprensa@otis.com;navigati_imprensa@otis.com;ouvidoria@otis.com'
$ws.Range('A10').Value = 'https://www.primac.com.br/manutencao-de-elevadores-sp.php'
$ws.Range('B10').Value = 'comercial@primac.com.br'
$ws.Range('A11').Value = 'https://www.rayteckelevadores.com.br/'
$ws.Range('B11').Value = 'rayteck@rayteckelevadores.com.br'
$ws.Range('A12').Value = 'https://nextelevadores.com/'
$ws.Range('B12').Value = 'hudsonkanegae@gatecubetecnologia.com'
$ws.Range('A13').Value = 'https://www.monciel.com.br/empresa-de-elevadores.php'
$ws.Range('B13').Value = 'monciel@monciel.com.br'
$ws.Range('A14').Value = 'https://crel.com.br/'
$ws.Range('B14').Value = 'bruno@crel.com.br'
$ws.Range('A15').Value = 'https://www.designelevadores.com.br/empresa-de-elevadores-em-sao-paulo.php'
$ws.Range('B15').Value = 'contato@designelevadores.com.br'
$ws.Range('A16').Value = 'https://www.elevadoreskorman.com.br/empresas-elevadores-sp'
$ws.Range('B16').Value = 'korman@elevadoreskorman.com.br;vendas@elevadoreskorman.com.br;comercial@elevadoreskorman.com.br'
$ws.Range('A17').Value = 'https://europaelevadores.com.br/manutencao-de-elevadores-em-sp/'
$ws.Range('B17').Value = 'europa@europaelevadores.com.br'
$ws.Range('A18').Value = 'https://www.surmonter.com.br/empresa-elevadores-sp'
$ws.Range('B18').Value = 'vendas@surmonter.com.br'
$ws.Range('A19').Value = 'https://www.ultronelevadores.com.br/empresas-elevadores-sp'
$ws.Range('B19').Value = 'contato@ultronelevadores.com.br;naoinformado@naoinformado.com'
$ws.Range('A20').Value = 'https://www.hts.com.br/'
$ws.Range('B20').Value = 'contato@hts.com.br'
$ws.Range('A21').Value = 'https://www.framartelelevadores.com.br/'
$ws.Range('B21').Value = 'framartelelevadores@terra.com.br;elcio_30@hotmail.com'
$ws.Range('A22').Value = 'https://www.mmelevadores.com.br/'
$ws.Range('B22').Value = 'orcamentos2@g7elevadores.com.br;comercial2@mmelevadores.com.br'
$ws.Range('A23').Value = 'https://www.gmvelevadores.com.br/hs/elevadores-em-sao-paulo/'
$ws.Range('B23').Value = 'contato@gmvelevadores.com.br'
$ws.Range('A24').Value = 'https://www.tecnewelevadores.com.br/manutencao-de-elevadores-em-sao-paulo'
$ws.Range('B24').Value = 'tecnica@tecnewelevadores.com.br'
$ws.Range('A25').Value = 'https://www.arsenalelevadores.com.br/empresas-elevadores-sp'
$ws.Range('B25').Value = 'arsenal@arsenalelevadores.com.br'
$ws.Range('A26').Value = 'https://espel.com.br/'
$ws.Range('B26').Value = 'espel@espel.com.br'
$ws.Range('A27').Value = 'https://atselevadores.com.br/'
$ws.Range('B27').Value = 'contato@atselevadores.com.br;info@atselevadores.com.br'
$ws.Range('A28').Value = 'https://www.flexst.com.br/empresa-elevadores-escadas-rolantes-sp'
$ws.Range('B28').Value = 'vendas@flexelevadores.com.br'
$ws.Range('A29').Value = 'https://novitaelevadores.com.br/'
$ws.Range('B29').Value = 'contato@novitaelevadores.com.br;novita@novitaelevadores.com.br'
$ws.Range('A30').Value = 'https://rcelevadores.com.br/'
$ws.Range('B30').Value = 'rcelevadores@hotmail.com'
$ws.Range('A31').Value = 'http://orionlift.com.br/'
$ws.Range('B31').Value = 'faleconosco@orionlift.com.br'
$ws.Range('A32').Value = 'https://elevatis.com.br/'
$ws.Range('B32').Value = 'rogerio.teodoro@elevatis.com.br'
$ws.Range('A33').Value = 'https://retrofitelevadores.com.br/'
$ws.Range('B33').Value = 'contato@elevadoresretrofit.com.br'
$ws.Range('A34').Value = 'https://www.astroselevadores.com.br/'
$ws.Range('B34').Value = 'contato@astroselevadores.com.br'
$ws.Range('A35').Value = 'https://sselev.com.br/'
$ws.Range('B35').Value = 'contato@sselev.com.br'
$ws.Range('A36').Value = 'https://www.elevadorestakaoki.com.br/'
$ws.Range('B36').Value = 'atendimento@elevadorestakaoki.com.br'
